# Updates canine keywords in the CypherOutput sheet, refreshes the
# StatOutput counters, adds a new sample/case (Glioma) row, clears the
# CaseDetailStat header row (files obj no longer has headers) and
# records that the files Cypher query was empty in CaseDetailStat_Message.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

# ---------------------------------------------------------------------
# CypherOutput (sheet1) : updated canine keywords + new Glioma case row
# ---------------------------------------------------------------------
$wsCypher = $wb.Worksheets.Item("CypherOutput")

$cypherRows = @(
    @("COTC007B-0503", "COTC007B", "Clinical Trial", "Beagle", "Lymphoma", "IIIa", "11.5", "Female", "Yes"),
    @("COTC007B-0211", "COTC007B", "Clinical Trial", "Beagle", "Lymphoma", "III", "9.9", "Male", "Yes"),
    @("COTC007B-0510", "COTC007B", "Clinical Trial", "Beagle", "Lymphoma", "IIIa", "9.5", "Male", "Yes"),
    @("COTC007B-0608", "COTC007B", "Clinical Trial", "Beagle", "Lymphoma", "IVa", "4.2", "Male", "Yes"),
    @("NCATS-COP01-CCB010045", "NCATS-COP01", "Transcriptomics", "Beagle", "T Cell Lymphoma", "III", "5.0", "Female", "Yes"),
    @("GLIOMA01-i_6561", "GLIOMA01", "Genomics", "Beagle", "Glioma", "Unknown", "8.0", "Male", "Yes")
)

$rowIndex = 2
foreach ($row in $cypherRows) {
    $colIndex = 1
    foreach ($val in $row) {
        $cell = $wsCypher.Cells.Item($rowIndex, $colIndex)
        if ($colIndex -eq 7) {
            # Age column holds numeric-looking text (e.g. "11.5") - keep it text
            Set-TextValue $cell $val
        } else {
            $cell.Value = $val
        }
        $colIndex++
    }
    $rowIndex++
}

# ---------------------------------------------------------------------
# StatOutput (sheet4) : refreshed counts
# ---------------------------------------------------------------------
$wsStat = $wb.Worksheets.Item("StatOutput")
Set-TextValue $wsStat.Cells.Item(2, 1) "12"
Set-TextValue $wsStat.Cells.Item(2, 2) "8"
Set-TextValue $wsStat.Cells.Item(2, 3) "6"
Set-TextValue $wsStat.Cells.Item(2, 4) "3"

# ---------------------------------------------------------------------
# CaseDetailStat (sheet6) : header row removed (files obj no headers)
# ---------------------------------------------------------------------
$wsCaseDetail = $wb.Worksheets.Item("CaseDetailStat")
$wsCaseDetail.Range("A1:F1").ClearContents()

# ---------------------------------------------------------------------
# CaseDetailStat_Message (sheet7) : files query was empty -> validation
# message logged, followed by the usual connection/cypher/output block
# with an empty Cypher string, shifting everything down by one row.
# ---------------------------------------------------------------------
$wsCaseMsg = $wb.Worksheets.Item("CaseDetailStat_Message")
$wsCaseMsg.Cells.Item(31, 1).Value = $wsCaseMsg.Cells.Item(30, 1).Value2
$wsCaseMsg.Cells.Item(30, 1).Value = $wsCaseMsg.Cells.Item(29, 1).Value2
$wsCaseMsg.Cells.Item(29, 1).Value = $wsCaseMsg.Cells.Item(28, 1).Value2
$wsCaseMsg.Cells.Item(28, 1).Value = $wsCaseMsg.Cells.Item(27, 1).Value2
$wsCaseMsg.Cells.Item(27, 1).Value = $wsCaseMsg.Cells.Item(26, 1).Value2
$wsCaseMsg.Cells.Item(26, 1).Value = $wsCaseMsg.Cells.Item(25, 1).Value2
$wsCaseMsg.Cells.Item(25, 1).Value = $wsCaseMsg.Cells.Item(24, 1).Value2
$wsCaseMsg.Cells.Item(24, 1).Value = $wsCaseMsg.Cells.Item(23, 1).Value2
$wsCaseMsg.Cells.Item(23, 1).Value = $wsCaseMsg.Cells.Item(22, 1).Value2
$wsCaseMsg.Cells.Item(22, 1).Value = $wsCaseMsg.Cells.Item(21, 1).Value2
$wsCaseMsg.Cells.Item(21, 1).Value = "Cypher query should not be an empty string"
$wsCaseMsg.Cells.Item(29, 1).Value = ""

Write-Host "Edit complete"
